{"js": "// Correct the appointments table: the FRANJA HORARIA (time slot) values\n// change, and the FLOR A FRUTO / INMERSSO BOUTIQUE appointments swap rows\n// (FLOR A FRUTO now comes first). The MESA column stays empty and the\n// third row's COMPRADOR (ENCADENAMIENTOS PRODUCTIVOS ...) is unchanged\n// apart from its time slot.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Data rows start at index 1 (row 0 is the FRANJA HORARIA/MESA/COMPRADOR\n// header row). Column 0 = FRANJA HORARIA, column 2 = COMPRADOR.\nconst newRows = [\n  { row: 1, time: \"08:30 - 08:45\", buyer: \"FLOR A FRUTO\" },\n  { row: 2, time: \"09:45 - 10:00\", buyer: \"INMERSSO BOUTIQUE\" },\n  { row: 3, time: \"10:30 - 10:45\", buyer: null } // buyer unchanged for this row\n];\n\nfor (const entry of newRows) {\n  const timeCell = table.getCell(entry.row, 0);\n  timeCell.value = entry.time;\n  await context.sync();\n\n  if (entry.buyer !== null) {\n    const buyerCell = table.getCell(entry.row, 2);\n    buyerCell.value = entry.buyer;\n    await context.sync();\n  }\n}\n", "ps1": "# Correct the appointments table in the document:\n#   - FRANJA HORARIA (time slot) values are updated\n#   - The FLOR A FRUTO / INMERSSO BOUTIQUE appointments swap rows\n#     (FLOR A FRUTO now scheduled first)\n#   - MESA stays blank; the third row's COMPRADOR\n#     (ENCADENAMIENTOS PRODUCTIVOS ...) keeps its buyer, only its time changes\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Table.Cell(row, column) is 1-indexed and row 1 is the header row\n# (FRANJA HORARIA | MESA | COMPRADOR), so data rows start at row 2.\n$t.Cell(2, 1).Range.Text = \"08:30 - 08:45\"\n$t.Cell(2, 3).Range.Text = \"FLOR A FRUTO\"\n\n$t.Cell(3, 1).Range.Text = \"09:45 - 10:00\"\n$t.Cell(3, 3).Range.Text = \"INMERSSO BOUTIQUE\"\n\n$t.Cell(4, 1).Range.Text = \"10:30 - 10:45\"\n"}
